$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Max sort:" block (rows 3-5): tweak a couple of values ---
$ws.Range("C3").Value = 1000
$ws.Range("H3").Value = 1000
$ws.Range("B4").Value = 100000
$ws.Range("G4").Value = 100000

# --- Remove the old trailing rows for the "Quick sort:" and "Radix sort:" blocks ---
# (row 11 = old 100,000,000 line of "Quick sort:"; row 17 = old 100,000,000 line of "Radix sort:")
$ws.Range("A11:I11").Clear()
$ws.Range("A17:I17").Clear()

# --- New row 6: "Quick sort:" label + first data line (array size 100000 / time 1000) ---
$ws.Range("A6").Value = "Quick sort:"
$ws.Range("B6").Value = 100000
$ws.Range("C6").Value = 1000
$ws.Range("F6").Value = "Quick sort:"
$ws.Range("G6").Value = 100000
$ws.Range("H6").Value = 1000

# Row 7 loses its "Quick sort:" label (now on row 6) but keeps its numbers
$ws.Range("A7").Clear()
$ws.Range("F7").Clear()

# Row 8 picks up the same cell style as row 9/10 (font + alignment), value unchanged
$ws.Range("B9").Copy()
$ws.Range("B8").PasteSpecial(-4122)
$ws.Range("G9").Copy()
$ws.Range("G8").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Rows 9 & 10 bump up to the next order of magnitude
$ws.Range("B9").Value = 10000000
$ws.Range("G9").Value = 10000000
$ws.Range("B10").Value = 100000000
$ws.Range("G10").Value = 100000000

# --- New row 12: "Radix sort:" label + first data line (array size 100000 / time 1000) ---
$ws.Range("A12").Value = "Radix sort:"
$ws.Range("B12").Value = 100000
$ws.Range("C12").Value = 1000
$ws.Range("F12").Value = "Radix sort:"
$ws.Range("G12").Value = 100000
$ws.Range("H12").Value = 1000

# Row 13 loses its "Radix sort:" label (now on row 12) but keeps its numbers
$ws.Range("A13").Clear()
$ws.Range("F13").Clear()

# Row 14 picks up the same cell style as row 15/16 (font + alignment), value unchanged
$ws.Range("B15").Copy()
$ws.Range("B14").PasteSpecial(-4122)
$ws.Range("G15").Copy()
$ws.Range("G14").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Rows 15 & 16 bump up to the next order of magnitude
$ws.Range("B15").Value = 10000000
$ws.Range("G15").Value = 10000000
$ws.Range("B16").Value = 100000000
$ws.Range("G16").Value = 100000000

# Match the selection left behind in the saved workbook
$ws.Range("E22").Select() | Out-Null
